$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 0.78
$ws.Range("F6").Value = 0.8363636363636364
$ws.Range("G14").Value = 0.7098120300751876
$ws.Range("F16").Value = 0.5729323308270677
$ws.Range("G16").Value = 0.6962406015037594
$ws.Range("F17").Value = 0.5714285714285714
$ws.Range("G17").Value = 0.6977443609022556
$ws.Range("F18").Value = 0.8322222222222222
$ws.Range("G18").Value = 0.9044444444444444
$ws.Range("F19").Value = 0.8066666666666666
$ws.Range("G19").Value = 0.8844444444444445
$ws.Range("F20").Value = 0.7988611111111111
$ws.Range("F21").Value = 0.7955555555555556
$ws.Range("G21").Value = 0.8744444444444445
$ws.Range("G25").Value = 0.9630978260869557
$ws.Range("G29").Value = 0.78
$ws.Range("F30").Value = 0.8272727272727273
$ws.Range("G31").Value = 0.9545454545454546
$ws.Range("F33").Value = 0.8181818181818181
$ws.Range("F42").Value = 0.5833333333333333
$ws.Range("G43").Value = 0.6909090909090909
$ws.Range("F44").Value = 0.5681818181818181
$ws.Range("G44").Value = 0.6909090909090909
$ws.Range("F46").Value = 0.8789473684210526
$ws.Range("F47").Value = 0.8789473684210526
$ws.Range("F49").Value = 0.8631578947368421
$ws.Range("F50").Value = 0.7850746268656716
$ws.Range("F51").Value = 0.7462686567164178
$ws.Range("F52").Value = 0.7343283582089553
$ws.Range("G52").Value = 0.8746268656716418
$ws.Range("G53").Value = 0.8746268656716418
$ws.Range("F54").Value = 0.792
$ws.Range("G54").Value = 0.904
$ws.Range("F55").Value = 0.7626666666666666
$ws.Range("F56").Value = 0.76
$ws.Range("F57").Value = 0.7493333333333334
$ws.Range("F63").Value = 0.8412698412698413
$ws.Range("G64").Value = 0.9650793650793651
$ws.Range("F65").Value = 0.8412698412698413
$ws.Range("F70").Value = 0.7930885529157667
$ws.Range("F71").Value = 0.7697624190064795
$ws.Range("G71").Value = 0.8267818574514039
$ws.Range("F73").Value = 0.7650107991360691
$ws.Range("G73").Value = 0.8220302375809936
